# Update the adult TOD raw-score -> standard-score lookup tables on all six
# age-band sheets: each table gains a new leading row for raw score 0 and the
# rest of the standard-score column is recomputed through raw score 20.

$wb = $excel.ActiveWorkbook

# raw score column (A2:A22) is identical across all six sheets: 0 .. 20
$rawScores = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20)

# standard score column (B2:B22) per sheet, indexed 1..6 to match tab order
$ssBySheet = @{
    1 = @(51,52,54,55,56,58,59,60,62,64,65,67,69,71,74,76,79,83,88,96,111)
    2 = @(51,53,54,55,56,58,59,61,62,64,66,67,69,71,74,76,79,83,88,96,111)
    3 = @(52,54,55,56,57,59,60,62,63,65,66,68,70,72,75,77,80,84,89,96,111)
    4 = @(54,55,56,58,59,60,62,63,65,66,68,70,72,74,76,78,81,85,90,97,111)
    5 = @(57,58,59,60,61,63,64,65,67,68,70,72,74,76,78,80,83,87,91,98,111)
    6 = @(62,63,64,65,67,68,69,70,72,73,75,76,78,80,82,84,87,90,94,100,112)
}

for ($sheetIdx = 1; $sheetIdx -le 6; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $ssValues = $ssBySheet[$sheetIdx]

    for ($i = 0; $i -lt $rawScores.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $rawScores[$i]
        $ws.Cells.Item($row, 2).Value = $ssValues[$i]
    }
}
